$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update canonical URL and Date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/StructureDefinition/fr-uf-role"
$meta.Range("B8").Value = "2025-05-05T08:11:38+00:00"

# --- Elements sheet: update Binding Value Set URL ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Z6").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-uf-role-code-for-prescription"

# Row 5 (Extension.url element) has a Fixed Value (column R) equal to the
# same canonical URL as Metadata!B2 - keep it in sync with the new URL.
$elements.Range("R5").Value = "https://hl7.fr/ig/fhir/medication/StructureDefinition/fr-uf-role"
